# Add a "BOUNDARY" attack block in columns AS:AZ (8 columns), mirroring the
# layout/style of the other attack blocks (e.g. the "FGSM" block in AK1:AR1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: merged header label "BOUNDARY" over AS1:AZ1 ---------------------
$ws.Range("AS1:AZ1").Merge()
$ws.Range("AK1:AR1").Copy()
$ws.Range("AS1:AZ1").PasteSpecial(-4122)
$ws.Range("AS1").Value = "BOUNDARY"

# --- Row 2: epsilon values (kept as text, matching the existing rows) ------
$epsilons = @("0.01", "0.02", "0.03", "0.04", "0.05", "0.07", "0.10", "0.20")
$ws.Range("AS2:AZ2").NumberFormat = "@"
for ($i = 0; $i -lt 8; $i++) {
    $col = 45 + $i   # AS=45 ... AZ=52
    $ws.Cells.Item(2, $col).Value = $epsilons[$i]
}
$ws.Range("AK2:AR2").Copy()
$ws.Range("AS2:AZ2").PasteSpecial(-4122)

# --- Data rows 4-12 (row 3 stays blank, like the other blocks) -------------
$data = @{
    4  = @(5.686472415924072, 5.738160610198975, 5.786428451538086, 5.824905395507812, 5.838682174682617, 6.180508613586426, 6.379770278930664, 8.032449722290039)
    5  = @(6.940323538590438, 6.971137078464893, 7.02943390174514, 7.062481095280284, 7.058417836854031, 7.450958640870901, 7.716292634911398, 10.01064954176535)
    6  = @(0.9995881319046021, 0.99958735704422, 0.9995734691619873, 0.9995725154876709, 0.9995645880699158, 0.9995023012161255, 0.9994622468948364, 0.9989767074584961)
    7  = @(5.008236408233643, 5.140683650970459, 5.160590648651123, 5.451739311218262, 5.69630241394043, 6.743438720703125, 8.162350654602051, 12.6733283996582)
    8  = @(5.947171656174707, 6.109762740453094, 6.155543872594171, 6.425935301629828, 6.895325180031679, 8.271070910932385, 10.15984799549205, 15.88342398063569)
    9  = @(0.9997783303260803, 0.9997398257255554, 0.9997555017471313, 0.9996582865715027, 0.9995759129524231, 0.9991956353187561, 0.9986026287078857, 0.9958498477935791)
    10 = @(2.663012027740479, 2.788084268569946, 3.055869340896606, 3.33517599105835, 3.796853542327881, 4.766287803649902, 5.873990058898926, 9.78786563873291)
    11 = @(3.583970138978879, 3.709712354300256, 3.937894407833777, 4.300472792426818, 4.849296298543166, 6.057199579677631, 7.353330425954725, 12.35969969350679)
    12 = @(0.9997842311859131, 0.999767005443573, 0.9997367262840271, 0.9996891021728516, 0.9995875954627991, 0.9993520379066467, 0.9990388154983521, 0.9973019957542419)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt 8; $i++) {
        $col = 45 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
